# This script reverts the multi-industry "AI/ML" themed training template
# back to its original generic "IT" themed content, and restores several
# blank separator rows that had been dropped near the top of each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: materialize a truly empty row (no cell data, no extra
# attributes) at the given 1-based row index on a worksheet, without
# shifting any existing rows. Toggling Hidden on/off forces the engine
# to write out the row element even though it carries no data.
# ---------------------------------------------------------------------
function Add-BlankRow {
    param($ws, [int]$rowIndex)
    $row = $ws.Rows.Item($rowIndex)
    $row.Hidden = $true
    $row.Hidden = $false
}

# ---------------------------------------------------------------------
# Sheet 1: "Training Schedule Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Schedule Overview")

# Row 3 was already a blank separator row in the source file; the COM
# round-trip silently drops untouched empty rows, so re-materialize it
# along with the two newly (re-)introduced blank rows.
Add-BlankRow $ws1 3
Add-BlankRow $ws1 6
Add-BlankRow $ws1 17

$ws1.Range("A9").Value = "IT Fundamentals (AI-101)"
$ws1.Range("A10").Value = "IT Platform Overview (AI-102)"
$ws1.Range("B11").Value = "System Administrators"
$ws1.Range("B12").Value = "IT Managers"
$ws1.Range("B13").Value = "DevOps Engineers, IT"
$ws1.Range("B14").Value = "DevOps Engineers, QA"

# ---------------------------------------------------------------------
# Sheet 2: "Detailed Training Schedule"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Detailed Training Schedule")

Add-BlankRow $ws2 2

$ws2.Range("B4").Value = "IT Fundamentals"
$ws2.Range("B5").Value = "IT Platform Overview"
$ws2.Range("C6").Value = "System Administrators"
$ws2.Range("C7").Value = "System Administrators"
$ws2.Range("C8").Value = "System Administrators"
$ws2.Range("C9").Value = "IT Managers"
$ws2.Range("C10").Value = "IT Managers"
$ws2.Range("C11").Value = "DevOps Engineers, IT"
$ws2.Range("C12").Value = "DevOps Engineers, IT"
$ws2.Range("C13").Value = "DevOps Engineers, QA"
$ws2.Range("C14").Value = "DevOps Engineers, QA"

# ---------------------------------------------------------------------
# Sheet 3: "Instructor Schedule"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Instructor Schedule")
Add-BlankRow $ws3 2

# ---------------------------------------------------------------------
# Sheet 4: "Facility Schedule"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Facility Schedule")
Add-BlankRow $ws4 2

# ---------------------------------------------------------------------
# Sheet 5: "Participant Tracking"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Participant Tracking")
Add-BlankRow $ws5 2

Write-Host "Edit complete"
